$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a text value (mirrors the original inlineStr
    # cells) even when the string looks numeric (e.g. "241.88"), without
    # leaving a stray style behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "96.895.68"
$ws.Range("E2").Value = "  +4.12%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.135.89"
$ws.Range("E3").Value = "  +0.46%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - Solana
Set-TextValue $ws.Range("D5") "241.88"
$ws.Range("E5").Value = "  +1.80%  "

# Row 6 - BNB
Set-TextValue $ws.Range("D6") "611.77"
$ws.Range("E6").Value = "  -0.33%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "1.12"
$ws.Range("E7").Value = "  +2.93%  "

# Row 8 - Dogecoin
Set-TextValue $ws.Range("D8") "0.384"
$ws.Range("E8").Value = "  -1.96%  "

# Row 10 - LidoStakedEther
Set-TextValue $ws.Range("D10") "3.131.58"
$ws.Range("E10").Value = "  +0.45%  "

# Row 11 - Cardano
Set-TextValue $ws.Range("D11") "0.781"
$ws.Range("E11").Value = "  -2.62%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.28%  "

# Row 13 - WrappedBTC
Set-TextValue $ws.Range("D13") "96.683.04"
$ws.Range("E13").Value = "  +4.27%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -1.89%  "

# Row 15 - was Avalanche, becomes Toncoin
$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D15") "5.56"
$ws.Range("E15").Value = "  +2.09%  "

# Row 16 - was Toncoin, becomes Avalanche
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D16") "34.13"
$ws.Range("E16").Value = "  -0.74%  "

# Row 17 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D17") "3.716.62"
$ws.Range("E17").Value = "  +0.21%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "3.135.99"
$ws.Range("E18").Value = "  +0.75%  "

# Row 19 - was BitcoinCash, becomes SuiNetwork
$ws.Range("B19").Value = "SuiNetwork"
$ws.Range("C19").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D19") "3.55"
$ws.Range("E19").Value = "  -6.27%  "

# Row 20 - was SuiNetwork, becomes BitcoinCash
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D20") "517.40"
$ws.Range("E20").Value = "  +17.32%  "

# Row 21 - Chainlink
Set-TextValue $ws.Range("D21") "14.59"
$ws.Range("E21").Value = "  +0.32%  "

# Row 22 - Polkadot
Set-TextValue $ws.Range("D22") "5.69"
$ws.Range("E22").Value = "  -3.84%  "

# Row 23 - PEPE
Set-TextValue $ws.Range("D23") "0.0000194"
$ws.Range("E23").Value = "  -4.64%  "

# Row 24 - Uniswap
Set-TextValue $ws.Range("D24") "8.84"
$ws.Range("E24").Value = "  -3.44%  "

# Row 25 - was Litecoin, becomes NEARProtocol
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D25") "5.49"
$ws.Range("E25").Value = "  -2.23%  "

# Row 26 - was NEARProtocol, becomes Litecoin
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D26") "88.73"
$ws.Range("E26").Value = "  +3.54%  "

# Row 27 - Aptos
Set-TextValue $ws.Range("D27") "11.62"
$ws.Range("E27").Value = "  -4.32%  "

# Row 28 - WrappedeETH
Set-TextValue $ws.Range("D28") "3.301.16"
$ws.Range("E28").Value = "  +0.42%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.09%  "

# Row 30 - Stellar
Set-TextValue $ws.Range("D30") "0.242"
$ws.Range("E30").Value = "  +2.87%  "

# Row 31 - Cronos
$ws.Range("E31").Value = "  -2.98%  "

# Row 32 - Hedera
Set-TextValue $ws.Range("D32") "0.126"
$ws.Range("E32").Value = "  +0.04%  "

# Row 33 - was InternetComputer(DFINITY), becomes Binance-PegBSC-USD
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D33") "0.997"
$ws.Range("E33").Value = "  -0.57%  "

# Row 34 - was Binance-PegBSC-USD, becomes InternetComputer(DFINITY)
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D34") "9.01"
$ws.Range("E34").Value = "  -1.89%  "

# Row 35 - EthereumClassic
Set-TextValue $ws.Range("D35") "26.68"
$ws.Range("E35").Value = "  +3.45%  "

# Row 36 - Kaspa
Set-TextValue $ws.Range("D36") "0.153"
$ws.Range("E36").Value = "  -4.48%  "

# Row 37 - RenderToken
Set-TextValue $ws.Range("D37") "7.34"
$ws.Range("E37").Value = "  -9.20%  "

# Row 38 - PancakeSwap
$ws.Range("E38").Value = "  -0.91%  "

# Row 39 - was Bittensor, becomes WhiteBITCoin
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D39") "24.19"
$ws.Range("E39").Value = "  +0.95%  "

# Row 40 - was WhiteBITCoin, becomes Bittensor
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D40") "478.55"
$ws.Range("E40").Value = "  +1.37%  "

# Row 41 - PolygonEcosystemToken
Set-TextValue $ws.Range("D41") "0.438"
$ws.Range("E41").Value = "  +1.75%  "

# Row 42 - Fetch.AI
Set-TextValue $ws.Range("D42") "1.22"
$ws.Range("E42").Value = "  -5.39%  "

# Row 43 - MantraDAO
Set-TextValue $ws.Range("D43") "3.58"
$ws.Range("E43").Value = "  -10.38%  "

# Row 44 - USDe
$ws.Range("E44").Value = "  -0.02%  "

# Row 45 - dogwifhat
Set-TextValue $ws.Range("D45") "3.18"
$ws.Range("E45").Value = "  -5.10%  "

# Row 46 - Monero
Set-TextValue $ws.Range("D46") "161.11"
$ws.Range("E46").Value = "  +1.83%  "

# Row 47 - Stacks
Set-TextValue $ws.Range("D47") "1.93"
$ws.Range("E47").Value = "  +4.56%  "

# Row 48 - ARBITRUM
Set-TextValue $ws.Range("D48") "0.703"
$ws.Range("E48").Value = "  +2.13%  "

# Row 49 - Filecoin
Set-TextValue $ws.Range("D49") "4.49"
$ws.Range("E49").Value = "  +2.97%  "

# Row 50 - OKB
Set-TextValue $ws.Range("D50") "44.32"
$ws.Range("E50").Value = "  +0.44%  "

# Row 51 - was VeChain, becomes FirstDigitalUSD
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D51") "0.998"
$ws.Range("E51").Value = "  +0.02%  "
